# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (fund-holding detail) right after the
# "总计" (totals) sheet, pushing "2022-Q2" / "2021-Q3" / "2021-Q1" back by
# one slot, and adds a matching summary row at the top of "总计".

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet    = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet right after "总计", using the existing
#    "2022-Q2" detail sheet as a formatting template (header style s="2",
#    index-column style s="2", thin borders, etc.) so the new sheet reuses
#    the same styles instead of minting new ones.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

$q2Sheet.Range("A1:H3").Copy($newSheet.Range("A1"))

# Row 4 needs to exist too (2022-Q3 has three data rows); clone row 3's
# formatting (index-column style + plain data cells) down one more row.
$q2Sheet.Range("A3:H3").Copy($newSheet.Range("A4"))

# Columns that hold numeric-looking values but must stay TEXT (fund codes
# with leading zeros, and the percentage/NAV figures which are text in the
# source data) — force Text format before writing so Excel doesn't coerce
# them into numbers.
$textCols = @("B", "D", "E", "F", "G")
foreach ($col in $textCols) {
    $newSheet.Range("$col" + "2:" + "$col" + "4").NumberFormat = "@"
}

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "166109"
$newSheet.Range("C2").Value = "信澳量化先锋混合（LOF）A"
$newSheet.Range("D2").Value = "0.79"
$newSheet.Range("E2").Value = "88.99"
$newSheet.Range("F2").Value = "2.37"
$newSheet.Range("G2").Value = "0.0187"
$newSheet.Range("H2").Value = 8

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "000398"
$newSheet.Range("C3").Value = "华富灵活配置混合"
$newSheet.Range("D3").Value = "0.12"
$newSheet.Range("E3").Value = "94.04"
$newSheet.Range("F3").Value = "3.26"
$newSheet.Range("G3").Value = "0.0039"
$newSheet.Range("H3").Value = 6

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "166110"
$newSheet.Range("C4").Value = "信澳量化先锋混合（LOF）C"
$newSheet.Range("D4").Value = "0.11"
$newSheet.Range("E4").Value = "88.99"
$newSheet.Range("F4").Value = "2.37"
$newSheet.Range("G4").Value = "0.0026"
$newSheet.Range("H4").Value = 8

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new top data row for 2022-Q3 and
#    shift the existing quarters down, re-numbering the A index column.
# ---------------------------------------------------------------------
# Extend the styled index-column (A) formatting down to the new row 5 by
# cloning row 4's cell.
$totalSheet.Range("A4").Copy($totalSheet.Range("A5"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.03

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.18

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.04

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q1"
$totalSheet.Range("C5").Value = 1
$totalSheet.Range("D5").Value = 0.02

$totalSheet.Range("A1").Select()
